$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 2: split off the trailing ". " into its own run "." ---
$para2 = $tr.Paragraphs(2, 1)
$tail = $para2.Characters($para2.Length - 1, 2)
$tail.Text = ""
$dot = $tr.InsertAfter(".")

# --- New paragraph 3: new sentence about the adaptive cruise control characteristics ---
$run1 = $tr.InsertAfter("`rZapoznaliśmy się z charakterystykami odległościowymi, prędkościowymi oraz ")
$run2 = $tr.InsertAfter("przyspieszeniowymi")
$run3 = $tr.InsertAfter(" ")
$run4 = $tr.InsertAfter("adaptacyjnego tempomatu.")

# --- The textbox has spAutoFit; PowerPoint grows it to fit the extra paragraph ---
$shp.Height = 236.28511811023623
